# Generate Report for Handoff
# Refresh the "Latest Handoff/Handback" timestamps for every row that is
# still pending (status "Ready for handoff") or errored ("Handback
# transform failed") on the Overview sheet and each locale sheet
# (zh-cn, de-de) - rows 7 and 10-16.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-21-18 12:21:44"
}

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "2016-03-18 12:21:37"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "2016-03-18 12:21:44"
}
